# Regenerate the "K" column (column G) values using the new calc (s_vals),
# replacing the previous "Strike#" derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
